$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Mdk"
$ws.Cells.Item(2, 3).Value2 = "Tspan1"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 2
$ws.Cells.Item(2, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(2, 7).Value2 = 0.437724
$ws.Cells.Item(2, 8).Value2 = 1.313172
$ws.Cells.Item(2, 9).Value2 = 0.02046276855287852
$ws.Cells.Item(2, 10).Value2 = 0.02204588088728605
$ws.Cells.Item(2, 11).Value2 = 1
$ws.Cells.Item(2, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 13).Value2 = 0.1180343333333333
$ws.Cells.Item(2, 14).Value2 = 0.354103
$ws.Cells.Item(2, 15).Value2 = 0.07902433777228687
$ws.Cells.Item(2, 16).Value2 = 0.08705977225985931
$ws.Cells.Item(2, 17).Value2 = 0.051666460524
$ws.Cells.Item(2, 18).Value2 = 0.464998144716
$ws.Cells.Item(2, 19).Value2 = 0.001617056733878802
$ws.Cells.Item(2, 20).Value2 = 0.001919309369315108

# Row 3
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Mdk"
$ws.Cells.Item(3, 3).Value2 = "Tspan1"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 2
$ws.Cells.Item(3, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(3, 7).Value2 = 0.437724
$ws.Cells.Item(3, 8).Value2 = 1.313172
$ws.Cells.Item(3, 9).Value2 = 0.02046276855287852
$ws.Cells.Item(3, 10).Value2 = 0.02204588088728605
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 0.8465543333333333
$ws.Cells.Item(3, 14).Value2 = 2.539663
$ws.Cells.Item(3, 15).Value2 = 0.5667706479182028
$ws.Cells.Item(3, 16).Value2 = 0.6244016074328403
$ws.Cells.Item(3, 17).Value2 = 0.370557149004
$ws.Cells.Item(3, 18).Value2 = 3.335014341036
$ws.Cells.Item(3, 19).Value2 = 0.01159769659091518
$ws.Cells.Item(3, 20).Value2 = 0.01376548346329434

# Row 4
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Mdk"
$ws.Cells.Item(4, 3).Value2 = "Tspan1"
$ws.Cells.Item(4, 4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value2 = 2
$ws.Cells.Item(4, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(4, 7).Value2 = 0.437724
$ws.Cells.Item(4, 8).Value2 = 1.313172
$ws.Cells.Item(4, 9).Value2 = 0.02046276855287852
$ws.Cells.Item(4, 10).Value2 = 0.02204588088728605
$ws.Cells.Item(4, 11).Value2 = 1
$ws.Cells.Item(4, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(4, 13).Value2 = 0.1154756666666667
$ws.Cells.Item(4, 14).Value2 = 0.346427
$ws.Cells.Item(4, 15).Value2 = 0.07731130281708999
$ws.Cells.Item(4, 16).Value2 = 0.08517255071170332
$ws.Cells.Item(4, 17).Value2 = 0.050546470716
$ws.Cells.Item(4, 18).Value2 = 0.454918236444
$ws.Cells.Item(4, 19).Value2 = 0.001582003296067618
$ws.Cells.Item(4, 20).Value2 = 0.001877703907856542

# Row 5
$ws.Cells.Item(5, 1).Value2 = "ECs"
$ws.Cells.Item(5, 2).Value2 = "Mdk"
$ws.Cells.Item(5, 3).Value2 = "Tspan1"
$ws.Cells.Item(5, 4).Value2 = "MuSCs"
$ws.Cells.Item(5, 5).Value2 = 2
$ws.Cells.Item(5, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(5, 7).Value2 = 0.437724
$ws.Cells.Item(5, 8).Value2 = 1.313172
$ws.Cells.Item(5, 9).Value2 = 0.02046276855287852
$ws.Cells.Item(5, 10).Value2 = 0.02204588088728605
$ws.Cells.Item(5, 11).Value2 = 1
$ws.Cells.Item(5, 12).Value2 = 0.5
$ws.Cells.Item(5, 13).Value2 = 0.413581
$ws.Cells.Item(5, 14).Value2 = 0.827162
$ws.Cells.Item(5, 15).Value2 = 0.2768937114924203
$ws.Cells.Item(5, 16).Value2 = 0.2033660695955972
$ws.Cells.Item(5, 17).Value2 = 0.181034329644
$ws.Cells.Item(5, 18).Value2 = 1.086205977864
$ws.Cells.Item(5, 19).Value2 = 0.005666011932016917
$ws.Cells.Item(5, 20).Value2 = 0.00448338414682006

# Row 6
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Mdk"
$ws.Cells.Item(6, 3).Value2 = "Tspan1"
$ws.Cells.Item(6, 4).Value2 = "ECs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 12.48419333333333
$ws.Cells.Item(6, 8).Value2 = 37.45258
$ws.Cells.Item(6, 9).Value2 = 0.5836124104444559
$ws.Cells.Item(6, 10).Value2 = 0.6287638767819841
$ws.Cells.Item(6, 11).Value2 = 1
$ws.Cells.Item(6, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(6, 13).Value2 = 0.1180343333333333
$ws.Cells.Item(6, 14).Value2 = 0.354103
$ws.Cells.Item(6, 15).Value2 = 0.07902433777228687
$ws.Cells.Item(6, 16).Value2 = 0.08705977225985931
$ws.Cells.Item(6, 17).Value2 = 1.473563437304444
$ws.Cells.Item(6, 18).Value2 = 13.26207093574
$ws.Cells.Item(6, 19).Value2 = 0.0461195842510612
$ws.Cells.Item(6, 20).Value2 = 0.05474003991786577

# Row 7
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Mdk"
$ws.Cells.Item(7, 3).Value2 = "Tspan1"
$ws.Cells.Item(7, 4).Value2 = "FAPs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 12.48419333333333
$ws.Cells.Item(7, 8).Value2 = 37.45258
$ws.Cells.Item(7, 9).Value2 = 0.5836124104444559
$ws.Cells.Item(7, 10).Value2 = 0.6287638767819841
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 0.8465543333333333
$ws.Cells.Item(7, 14).Value2 = 2.539663
$ws.Cells.Item(7, 15).Value2 = 0.5667706479182028
$ws.Cells.Item(7, 16).Value2 = 0.6244016074328403
$ws.Cells.Item(7, 17).Value2 = 10.56854796450444
$ws.Cells.Item(7, 18).Value2 = 95.11693168053999
$ws.Cells.Item(7, 19).Value2 = 0.3307743840007084
$ws.Cells.Item(7, 20).Value2 = 0.3926011753583752

# Row 8
$ws.Cells.Item(8, 1).Value2 = "FAPs"
$ws.Cells.Item(8, 2).Value2 = "Mdk"
$ws.Cells.Item(8, 3).Value2 = "Tspan1"
$ws.Cells.Item(8, 4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 12.48419333333333
$ws.Cells.Item(8, 8).Value2 = 37.45258
$ws.Cells.Item(8, 9).Value2 = 0.5836124104444559
$ws.Cells.Item(8, 10).Value2 = 0.6287638767819841
$ws.Cells.Item(8, 11).Value2 = 1
$ws.Cells.Item(8, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(8, 13).Value2 = 0.1154756666666667
$ws.Cells.Item(8, 14).Value2 = 0.346427
$ws.Cells.Item(8, 15).Value2 = 0.07731130281708999
$ws.Cells.Item(8, 16).Value2 = 0.08517255071170332
$ws.Cells.Item(8, 17).Value2 = 1.441620547962222
$ws.Cells.Item(8, 18).Value2 = 12.97458493166
$ws.Cells.Item(8, 19).Value2 = 0.04511983579168314
$ws.Cells.Item(8, 20).Value2 = 0.05355342318090072

# Row 9
$ws.Cells.Item(9, 1).Value2 = "FAPs"
$ws.Cells.Item(9, 2).Value2 = "Mdk"
$ws.Cells.Item(9, 3).Value2 = "Tspan1"
$ws.Cells.Item(9, 4).Value2 = "MuSCs"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 12.48419333333333
$ws.Cells.Item(9, 8).Value2 = 37.45258
$ws.Cells.Item(9, 9).Value2 = 0.5836124104444559
$ws.Cells.Item(9, 10).Value2 = 0.6287638767819841
$ws.Cells.Item(9, 11).Value2 = 1
$ws.Cells.Item(9, 12).Value2 = 0.5
$ws.Cells.Item(9, 13).Value2 = 0.413581
$ws.Cells.Item(9, 14).Value2 = 0.827162
$ws.Cells.Item(9, 15).Value2 = 0.2768937114924203
$ws.Cells.Item(9, 16).Value2 = 0.2033660695955972
$ws.Cells.Item(9, 17).Value2 = 5.163225162993332
$ws.Cells.Item(9, 18).Value2 = 30.97935097796
$ws.Cells.Item(9, 19).Value2 = 0.1615986064010032
$ws.Cells.Item(9, 20).Value2 = 0.1278692383248425

# Row 10
$ws.Cells.Item(10, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value2 = "Mdk"
$ws.Cells.Item(10, 3).Value2 = "Tspan1"
$ws.Cells.Item(10, 4).Value2 = "ECs"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 1.796802333333333
$ws.Cells.Item(10, 8).Value2 = 5.390407
$ws.Cells.Item(10, 9).Value2 = 0.08399710841140098
$ws.Cells.Item(10, 10).Value2 = 0.09049558675938332
$ws.Cells.Item(10, 11).Value2 = 1
$ws.Cells.Item(10, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(10, 13).Value2 = 0.1180343333333333
$ws.Cells.Item(10, 14).Value2 = 0.354103
$ws.Cells.Item(10, 15).Value2 = 0.07902433777228687
$ws.Cells.Item(10, 16).Value2 = 0.08705977225985931
$ws.Cells.Item(10, 17).Value2 = 0.2120843655467778
$ws.Cells.Item(10, 18).Value2 = 1.908759289921
$ws.Cells.Item(10, 19).Value2 = 0.00663781586699795
$ws.Cells.Item(10, 20).Value2 = 0.007878525173794252

# Row 11
$ws.Cells.Item(11, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(11, 2).Value2 = "Mdk"
$ws.Cells.Item(11, 3).Value2 = "Tspan1"
$ws.Cells.Item(11, 4).Value2 = "FAPs"
$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 6).Value2 = 1
$ws.Cells.Item(11, 7).Value2 = 1.796802333333333
$ws.Cells.Item(11, 8).Value2 = 5.390407
$ws.Cells.Item(11, 9).Value2 = 0.08399710841140098
$ws.Cells.Item(11, 10).Value2 = 0.09049558675938332
$ws.Cells.Item(11, 11).Value2 = 3
$ws.Cells.Item(11, 12).Value2 = 1
$ws.Cells.Item(11, 13).Value2 = 0.8465543333333333
$ws.Cells.Item(11, 14).Value2 = 2.539663
$ws.Cells.Item(11, 15).Value2 = 0.5667706479182028
$ws.Cells.Item(11, 16).Value2 = 0.6244016074328403
$ws.Cells.Item(11, 17).Value2 = 1.521090801426778
$ws.Cells.Item(11, 18).Value2 = 13.689817212841
$ws.Cells.Item(11, 19).Value2 = 0.04760709555758525
$ws.Cells.Item(11, 20).Value2 = 0.05650558983813701

# Row 12
$ws.Cells.Item(12, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value2 = "Mdk"
$ws.Cells.Item(12, 3).Value2 = "Tspan1"
$ws.Cells.Item(12, 4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 6).Value2 = 1
$ws.Cells.Item(12, 7).Value2 = 1.796802333333333
$ws.Cells.Item(12, 8).Value2 = 5.390407
$ws.Cells.Item(12, 9).Value2 = 0.08399710841140098
$ws.Cells.Item(12, 10).Value2 = 0.09049558675938332
$ws.Cells.Item(12, 11).Value2 = 1
$ws.Cells.Item(12, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(12, 13).Value2 = 0.1154756666666667
$ws.Cells.Item(12, 14).Value2 = 0.346427
$ws.Cells.Item(12, 15).Value2 = 0.07731130281708999
$ws.Cells.Item(12, 16).Value2 = 0.08517255071170332
$ws.Cells.Item(12, 17).Value2 = 0.2074869473098888
$ws.Cells.Item(12, 18).Value2 = 1.867382525789
$ws.Cells.Item(12, 19).Value2 = 0.006493925884153758
$ws.Cells.Item(12, 20).Value2 = 0.007707739952448923

# Row 13
$ws.Cells.Item(13, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value2 = "Mdk"
$ws.Cells.Item(13, 3).Value2 = "Tspan1"
$ws.Cells.Item(13, 4).Value2 = "MuSCs"
$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 6).Value2 = 1
$ws.Cells.Item(13, 7).Value2 = 1.796802333333333
$ws.Cells.Item(13, 8).Value2 = 5.390407
$ws.Cells.Item(13, 9).Value2 = 0.08399710841140098
$ws.Cells.Item(13, 10).Value2 = 0.09049558675938332
$ws.Cells.Item(13, 11).Value2 = 1
$ws.Cells.Item(13, 12).Value2 = 0.5
$ws.Cells.Item(13, 13).Value2 = 0.413581
$ws.Cells.Item(13, 14).Value2 = 0.827162
$ws.Cells.Item(13, 15).Value2 = 0.2768937114924203
$ws.Cells.Item(13, 16).Value2 = 0.2033660695955972
$ws.Cells.Item(13, 17).Value2 = 0.7431233058223332
$ws.Cells.Item(13, 18).Value2 = 4.458739834934
$ws.Cells.Item(13, 19).Value2 = 0.02325827110266402
$ws.Cells.Item(13, 20).Value2 = 0.01840373179500315

# Row 14
$ws.Cells.Item(14, 1).Value2 = "MuSCs"
$ws.Cells.Item(14, 2).Value2 = "Mdk"
$ws.Cells.Item(14, 3).Value2 = "Tspan1"
$ws.Cells.Item(14, 4).Value2 = "ECs"
$ws.Cells.Item(14, 5).Value2 = 2
$ws.Cells.Item(14, 6).Value2 = 1
$ws.Cells.Item(14, 7).Value2 = 4.608308
$ws.Cells.Item(14, 8).Value2 = 9.216616
$ws.Cells.Item(14, 9).Value2 = 0.2154296772038511
$ws.Cells.Item(14, 10).Value2 = 0.154731001361478
$ws.Cells.Item(14, 11).Value2 = 1
$ws.Cells.Item(14, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(14, 13).Value2 = 0.1180343333333333
$ws.Cells.Item(14, 14).Value2 = 0.354103
$ws.Cells.Item(14, 15).Value2 = 0.07902433777228687
$ws.Cells.Item(14, 16).Value2 = 0.08705977225985931
$ws.Cells.Item(14, 17).Value2 = 0.5439385625746667
$ws.Cells.Item(14, 18).Value2 = 3.263631375448
$ws.Cells.Item(14, 19).Value2 = 0.01702418757753186
$ws.Cells.Item(14, 20).Value2 = 0.01347084574007026

# Row 15
$ws.Cells.Item(15, 1).Value2 = "MuSCs"
$ws.Cells.Item(15, 2).Value2 = "Mdk"
$ws.Cells.Item(15, 3).Value2 = "Tspan1"
$ws.Cells.Item(15, 4).Value2 = "FAPs"
$ws.Cells.Item(15, 5).Value2 = 2
$ws.Cells.Item(15, 6).Value2 = 1
$ws.Cells.Item(15, 7).Value2 = 4.608308
$ws.Cells.Item(15, 8).Value2 = 9.216616
$ws.Cells.Item(15, 9).Value2 = 0.2154296772038511
$ws.Cells.Item(15, 10).Value2 = 0.154731001361478
$ws.Cells.Item(15, 11).Value2 = 3
$ws.Cells.Item(15, 12).Value2 = 1
$ws.Cells.Item(15, 13).Value2 = 0.8465543333333333
$ws.Cells.Item(15, 14).Value2 = 2.539663
$ws.Cells.Item(15, 15).Value2 = 0.5667706479182028
$ws.Cells.Item(15, 16).Value2 = 0.6244016074328403
$ws.Cells.Item(15, 17).Value2 = 3.901183106734667
$ws.Cells.Item(15, 18).Value2 = 23.407098640408
$ws.Cells.Item(15, 19).Value2 = 0.122099217729636
$ws.Cells.Item(15, 20).Value2 = 0.09661428596979987

# Row 16
$ws.Cells.Item(16, 1).Value2 = "MuSCs"
$ws.Cells.Item(16, 2).Value2 = "Mdk"
$ws.Cells.Item(16, 3).Value2 = "Tspan1"
$ws.Cells.Item(16, 4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value2 = 2
$ws.Cells.Item(16, 6).Value2 = 1
$ws.Cells.Item(16, 7).Value2 = 4.608308
$ws.Cells.Item(16, 8).Value2 = 9.216616
$ws.Cells.Item(16, 9).Value2 = 0.2154296772038511
$ws.Cells.Item(16, 10).Value2 = 0.154731001361478
$ws.Cells.Item(16, 11).Value2 = 1
$ws.Cells.Item(16, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(16, 13).Value2 = 0.1154756666666667
$ws.Cells.Item(16, 14).Value2 = 0.346427
$ws.Cells.Item(16, 15).Value2 = 0.07731130281708999
$ws.Cells.Item(16, 16).Value2 = 0.08517255071170332
$ws.Cells.Item(16, 17).Value2 = 0.5321474385053333
$ws.Cells.Item(16, 18).Value2 = 3.192884631032
$ws.Cells.Item(16, 19).Value2 = 0.01665514901009488
$ws.Cells.Item(16, 20).Value2 = 0.01317883406013312

# Row 17
$ws.Cells.Item(17, 1).Value2 = "MuSCs"
$ws.Cells.Item(17, 2).Value2 = "Mdk"
$ws.Cells.Item(17, 3).Value2 = "Tspan1"
$ws.Cells.Item(17, 4).Value2 = "MuSCs"
$ws.Cells.Item(17, 5).Value2 = 2
$ws.Cells.Item(17, 6).Value2 = 1
$ws.Cells.Item(17, 7).Value2 = 4.608308
$ws.Cells.Item(17, 8).Value2 = 9.216616
$ws.Cells.Item(17, 9).Value2 = 0.2154296772038511
$ws.Cells.Item(17, 10).Value2 = 0.154731001361478
$ws.Cells.Item(17, 11).Value2 = 1
$ws.Cells.Item(17, 12).Value2 = 0.5
$ws.Cells.Item(17, 13).Value2 = 0.413581
$ws.Cells.Item(17, 14).Value2 = 0.827162
$ws.Cells.Item(17, 15).Value2 = 0.2768937114924203
$ws.Cells.Item(17, 16).Value2 = 0.2033660695955972
$ws.Cells.Item(17, 17).Value2 = 1.905908630948
$ws.Cells.Item(17, 18).Value2 = 7.623634523792
$ws.Cells.Item(17, 19).Value2 = 0.05965112288658839
$ws.Cells.Item(17, 20).Value2 = 0.03146703559147478

# Row 18
$ws.Cells.Item(18, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(18, 2).Value2 = "Mdk"
$ws.Cells.Item(18, 3).Value2 = "Tspan1"
$ws.Cells.Item(18, 4).Value2 = "ECs"
$ws.Cells.Item(18, 5).Value2 = 2
$ws.Cells.Item(18, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(18, 7).Value2 = 2.064212666666667
$ws.Cells.Item(18, 8).Value2 = 6.192638000000001
$ws.Cells.Item(18, 9).Value2 = 0.09649803538741349
$ws.Cells.Item(18, 10).Value2 = 0.1039636542098684
$ws.Cells.Item(18, 11).Value2 = 1
$ws.Cells.Item(18, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(18, 13).Value2 = 0.1180343333333333
$ws.Cells.Item(18, 14).Value2 = 0.354103
$ws.Cells.Item(18, 15).Value2 = 0.07902433777228687
$ws.Cells.Item(18, 16).Value2 = 0.08705977225985931
$ws.Cells.Item(18, 17).Value2 = 0.2436479659682222
$ws.Cells.Item(18, 18).Value2 = 2.192831693714
$ws.Cells.Item(18, 19).Value2 = 0.007625693342817055
$ws.Cells.Item(18, 20).Value2 = 0.009051052058813907

# Row 19
$ws.Cells.Item(19, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(19, 2).Value2 = "Mdk"
$ws.Cells.Item(19, 3).Value2 = "Tspan1"
$ws.Cells.Item(19, 4).Value2 = "FAPs"
$ws.Cells.Item(19, 5).Value2 = 2
$ws.Cells.Item(19, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(19, 7).Value2 = 2.064212666666667
$ws.Cells.Item(19, 8).Value2 = 6.192638000000001
$ws.Cells.Item(19, 9).Value2 = 0.09649803538741349
$ws.Cells.Item(19, 10).Value2 = 0.1039636542098684
$ws.Cells.Item(19, 11).Value2 = 3
$ws.Cells.Item(19, 12).Value2 = 1
$ws.Cells.Item(19, 13).Value2 = 0.8465543333333333
$ws.Cells.Item(19, 14).Value2 = 2.539663
$ws.Cells.Item(19, 15).Value2 = 0.5667706479182028
$ws.Cells.Item(19, 16).Value2 = 0.6244016074328403
$ws.Cells.Item(19, 17).Value2 = 1.747468177888222
$ws.Cells.Item(19, 18).Value2 = 15.727213600994
$ws.Cells.Item(19, 19).Value2 = 0.054692254039358
$ws.Cells.Item(19, 20).Value2 = 0.0649150728032338

# Row 20
$ws.Cells.Item(20, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(20, 2).Value2 = "Mdk"
$ws.Cells.Item(20, 3).Value2 = "Tspan1"
$ws.Cells.Item(20, 4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(20, 5).Value2 = 2
$ws.Cells.Item(20, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(20, 7).Value2 = 2.064212666666667
$ws.Cells.Item(20, 8).Value2 = 6.192638000000001
$ws.Cells.Item(20, 9).Value2 = 0.09649803538741349
$ws.Cells.Item(20, 10).Value2 = 0.1039636542098684
$ws.Cells.Item(20, 11).Value2 = 1
$ws.Cells.Item(20, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(20, 13).Value2 = 0.1154756666666667
$ws.Cells.Item(20, 14).Value2 = 0.346427
$ws.Cells.Item(20, 15).Value2 = 0.07731130281708999
$ws.Cells.Item(20, 16).Value2 = 0.08517255071170332
$ws.Cells.Item(20, 17).Value2 = 0.2383663338251111
$ws.Cells.Item(20, 18).Value2 = 2.145297004426
$ws.Cells.Item(20, 19).Value2 = 0.00746038883509059
$ws.Cells.Item(20, 20).Value2 = 0.008854849610364006

# Row 21
$ws.Cells.Item(21, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(21, 2).Value2 = "Mdk"
$ws.Cells.Item(21, 3).Value2 = "Tspan1"
$ws.Cells.Item(21, 4).Value2 = "MuSCs"
$ws.Cells.Item(21, 5).Value2 = 2
$ws.Cells.Item(21, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(21, 7).Value2 = 2.064212666666667
$ws.Cells.Item(21, 8).Value2 = 6.192638000000001
$ws.Cells.Item(21, 9).Value2 = 0.09649803538741349
$ws.Cells.Item(21, 10).Value2 = 0.1039636542098684
$ws.Cells.Item(21, 11).Value2 = 1
$ws.Cells.Item(21, 12).Value2 = 0.5
$ws.Cells.Item(21, 13).Value2 = 0.413581
$ws.Cells.Item(21, 14).Value2 = 0.827162
$ws.Cells.Item(21, 15).Value2 = 0.2768937114924203
$ws.Cells.Item(21, 16).Value2 = 0.2033660695955972
$ws.Cells.Item(21, 17).Value2 = 0.8537191388926666
$ws.Cells.Item(21, 18).Value2 = 5.122314833356
$ws.Cells.Item(21, 19).Value2 = 0.02671969917014784
$ws.Cells.Item(21, 20).Value2 = 0.0211426797374567
